$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 1 ("TextBox 3"): merge the "Let’s " + "continue at " runs
#     into a single run, in the last paragraph of the text box.
$shp1 = $s.Shapes.Item(1)
$tr1 = $shp1.TextFrame.TextRange
$lastPara = $tr1.Paragraphs($tr1.Paragraphs().Count, 1)
$mergeRange = $tr1.Characters($lastPara.Start, 18)
$mergeRange.Text = "Let’s continue at "
# Editing the run structure nudges this autofit textbox's cached height;
# restore the original (unchanged-by-design) height explicitly.
$shp1.Height = 390.173480

# --- Shape 4 ("TextBox 8"): replace the header sentence with new wording
#     and shrink the text box to fit the shorter text.
$shp4 = $s.Shapes.Item(4)
$shp4.TextFrame.TextRange.Text = "Integration of web computing in scientific publishing"
$shp4.Width = 360.644730
